$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("obj")

# Update existing values and append new rows (A1:A10)
$values = @(87369811, 87369552, 87369391, 87369170, 87369030, 87368891, 87368629, 87145832, 87162117, 87364267)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
    # match the explicit 25.5pt custom row height used throughout the sheet
    $ws.Rows.Item($row).RowHeight = 25.5
}

# Enable iterative calculation with a delta of 1E-4 (matches workbook calcPr change)
$excel.Iteration = $true
$excel.MaxChange = 0.0001

# Update the active selection to match the post-edit state (A13)
$ws.Range("A13").Select()
